$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.5544143368907429
$ws.Range("J2").Value = 0.5544143368907429
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 5.833790839957222
$ws.Range("R2").Value = 52.504117559615
$ws.Range("S2").Value = 0.1784322478211398
$ws.Range("T2").Value = 0.1784322478211398

# Row 3
$ws.Range("I3").Value = 0.5544143368907429
$ws.Range("J3").Value = 0.5544143368907429
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("S3").Value = 0.2400054062599131
$ws.Range("T3").Value = 0.2400054062599131

# Row 4
$ws.Range("I4").Value = 0.5544143368907429
$ws.Range("J4").Value = 0.5544143368907429
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 4.445718396251444
$ws.Range("R4").Value = 40.011465566263
$ws.Range("S4").Value = 0.13597668280969
$ws.Range("T4").Value = 0.13597668280969

# Row 5
$ws.Range("G5").Value = 0.1587963333333333
$ws.Range("H5").Value = 0.476389
$ws.Range("I5").Value = 0.4455856631092571
$ws.Range("J5").Value = 0.4455856631092571
$ws.Range("M5").Value = 29.52617166666667
$ws.Range("N5").Value = 88.57851500000001
$ws.Range("O5").Value = 0.3218391660320701
$ws.Range("P5").Value = 0.3218391660320701
$ws.Range("Q5").Value = 4.688647798037223
$ws.Range("R5").Value = 42.197830182335
$ws.Range("S5").Value = 0.1434069182109302
$ws.Range("T5").Value = 0.1434069182109302

# Row 6
$ws.Range("G6").Value = 0.1587963333333333
$ws.Range("H6").Value = 0.476389
$ws.Range("I6").Value = 0.4455856631092571
$ws.Range("J6").Value = 0.4455856631092571
$ws.Range("O6").Value = 0.4328989896002822
$ws.Range("P6").Value = 0.4328989896002822
$ws.Range("Q6").Value = 6.306600030649001
$ws.Range("R6").Value = 56.759400275841
$ws.Range("S6").Value = 0.1928935833403691
$ws.Range("T6").Value = 0.1928935833403692

# Row 7
$ws.Range("G7").Value = 0.1587963333333333
$ws.Range("H7").Value = 0.476389
$ws.Range("I7").Value = 0.4455856631092571
$ws.Range("J7").Value = 0.4455856631092571
$ws.Range("M7").Value = 22.50081433333333
$ws.Range("N7").Value = 67.502443
$ws.Range("O7").Value = 0.2452618443676477
$ws.Range("P7").Value = 0.2452618443676476
$ws.Range("Q7").Value = 3.573046813147445
$ws.Range("R7").Value = 32.157421318327
$ws.Range("S7").Value = 0.1092851615579577
$ws.Range("T7").Value = 0.1092851615579577
